$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell text ---
$ws.Range("F7").Value = "string _analisingText"
$ws.Range("B10").Value = "pub StartUp(void)"
$ws.Range("F10").Value = "pub Train()"
$ws.Range("B11").Value = "pub Result(string mostViable, double percent)"

# --- New cell in new "FileObj" box (top-right, header only) ---
$ws.Range("J2").Value = "FileObj"
$ws.Range("J2").Font.Bold = $true

# --- New attribute under Menu class ---
$ws.Range("B6").Value = "BayesingNetwork _bn"

# --- New method under Menu class ---
$ws.Range("B12").Value = "pub SaveBayesingNetwork()"

# --- Clear old "FileReader" box cells (moved/rebuilt at column H) ---
$ws.Range("B14").Clear()
$ws.Range("B15").Clear()
$ws.Range("B16").Clear()

# --- Clear old loose method cells (moved into new FileReader box at column H) ---
$ws.Range("B18").Clear()
$ws.Range("B19").Clear()

# --- Rebuild "FileReader" class box at column H ---
$ws.Range("H5").Value = "FileReader"
$ws.Range("H5").Font.Bold = $true
$ws.Range("H6").Value = "string _BayesingNetworkFolder"
$ws.Range("H7").Value = "string _TestDataFolder"
$ws.Range("H8").Value = "string _TrainingDataFolder"
$ws.Range("H10").Value = "pub FileObj[] getTestData()"
$ws.Range("H11").Value = "pub FileObj[] GetSavedBayesingNetworks()"
$ws.Range("H12").Value = "pub FileObj[] GetTrainingData()"
$ws.Range("H13").Value = "pub string[] GetLemmatizingWords()"
$ws.Range("H14").Value = "pub bool SaveBayesingToFile(string folderName, List<CategoryObj> bayesingNetwork)"
$ws.Range("H15").Value = "pub BayesingNetwork[] GetSavedBayesingNetworks()"

# --- Column widths (resized to fit the new, longer class-diagram text) ---
$ws.Columns.Item(4).ColumnWidth = 43
$ws.Columns.Item(6).ColumnWidth = 32.85546875
$ws.Columns.Item(8).ColumnWidth = 75.7109375
$ws.Columns.Item(10).ColumnWidth = 7.42578125

# --- Update selection to reflect where the author left off ---
$ws.Range("F11").Select()
